$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 33, shifting existing rows 33:39 down to 34:40
$ws.Rows("33:33").Insert()

# Populate the newly inserted row 33 with the new record
$ws.Range("A33").Value = 5
$ws.Range("B33").Value = "Macroferia Regional de Talca"
$ws.Range("C33").Value = "Maule"
$ws.Range("D33").Value = 44776
$ws.Range("E33").Value = 7
$ws.Range("F33").Value = 100112040
$ws.Range("G33").Value = "Cilantro"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 150
$ws.Range("K33").Value = 10000
$ws.Range("L33").Value = 10000
$ws.Range("M33").Value = 10000
$ws.Range("N33").Value = '$/caja 36 atados'
$ws.Range("O33").Value = "Región de Coquimbo"
$ws.Range("P33").Value = 278
$ws.Range("Q33").Value = 36
$ws.Range("R33").Value = "Hortaliza"

# Match the date number format used by the rest of column D
$ws.Range("D33").NumberFormat = $ws.Range("D34").NumberFormat
